# Add predict weather modules
# Appends 14 new Speech/Command rows (278-291) to Sheet1, mirroring the
# formatting of the existing data rows (style + row height) and finishes
# with the selection left on the first empty row below the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Speech (col A) / Command label (col B) pairs to append, in order.
$rows = @(
    @("Tôi muốn buộc gió lại cho hương đừng bay đi", "undefined"),
    @("Sao nó không chạy", "undefined"),
    @("Con mèo đang vật lộn với quả len", "undefined"),
    @("turn air conditioner on", "bật điều hòa"),
    @("cooling the air", "bật điều hòa"),
    @("cool air", "bật điều hòa"),
    @("air conditioner on", "bật điều hòa"),
    @("turn air conditioner off", "tắt điều hòa"),
    @("off the air", "tắt điều hòa"),
    @("off the conditioner", "tắt điều hòa"),
    @("on the conditioner", "bật điều hòa"),
    @("too light", "tắt đèn chùm"),
    @("enough light for me assistant", "tắt đèn chùm"),
    @("Tôi thấy sáng như vậy là quá đủ hãy tắt đèn đi", "tắt đèn chùm")
)

# Last populated row before the insert.
$lastRow = 277
$row = $lastRow

foreach ($pair in $rows) {
    $row = $row + 1

    # Copy the formatting (font/style + row height) from the previous row
    # down onto the new row before writing values, so the new rows match
    # the look of the existing table (style index reused, ht="17.25").
    $ws.Range("A" + ($row - 1) + ":B" + ($row - 1)).Copy()
    $ws.Range("A" + $row + ":B" + $row).PasteSpecial(-4122)
    $ws.Rows($row).RowHeight = $ws.Rows($row - 1).RowHeight

    $ws.Range("A" + $row).Value = $pair[0]
    $ws.Range("B" + $row).Value = $pair[1]
}

$excel.CutCopyMode = $false

# Leave the selection / view like the saved workbook: first empty row
# right below the newly appended data.
$newRow = $row + 1
$ws.Range("A" + $newRow).Select()
